$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.530281782150269
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 1.999173641204834
$ws.Range("D1").Value = 1.221300601959229
$ws.Range("E1").Value = 0.9600273966789246
